# Captain Tsubasa py — "solo implemented club shared implemented new locator
# with imshow and many more ..." commit: drop the now-unused OCR sheet,
# refresh two difficulty-template descriptions (-> "... horizontal"), bump
# template 010's box coordinates, and append eight new UI-locator templates
# (031-038) for the club "shared play" flow.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$templates = $wb.Worksheets.Item("Templates")
$points = $wb.Worksheets.Item("Points")

# --- 1. Remove the obsolete "OCR" sheet -------------------------------------
$ocr = $wb.Worksheets.Item("OCR")
$ocr.Delete() | Out-Null

# --- 2. Templates sheet: tweak existing rows --------------------------------
# Template 010 bounding box changed.
$templates.Range("B4").Value = 1082
$templates.Range("C4").Value = 235
$templates.Range("D4").Value = 1410
$templates.Range("E4").Value = 288

# Clarify the two difficulty templates as "horizontal" variants now that a
# "vertical" counterpart exists (template 031 below).
$templates.Range("F8").Value = "difficulty -> very hard horizontal"
$templates.Range("F26").Value = "difficulty -> extreme horizontal"

# --- 3. Templates sheet: append the new club-shared-play locators ----------
$newRows = @(
    @("031", 696,  269, 810,  350, "difficulty -> very hard vertical "),
    @("032", 1550, 926, 1683, 1015, "club shared playe button"),
    @("033", 330,  187, 441,  260, "club shared play - accepting member list -> rank"),
    @("034", 1444, 735, 1585, 808, "join button"),
    @("035", 92,   66,  603,  109, "club shared play - accepting member list title"),
    @("036", 769,  327, 1089, 371, "failed to join dialog -> title"),
    @("037", 877,  677, 975,  724, "failed to join dialog -> ok button"),
    @("038", 767,  952, 1078, 999, "go to scenario list ->  shared play")
)

$row = 32
foreach ($r in $newRows) {
    $templates.Cells.Item($row, 1).Value = $r[0]
    $templates.Cells.Item($row, 2).Value = $r[1]
    $templates.Cells.Item($row, 3).Value = $r[2]
    $templates.Cells.Item($row, 4).Value = $r[3]
    $templates.Cells.Item($row, 5).Value = $r[4]
    $templates.Cells.Item($row, 6).Value = $r[5]
    $row = $row + 1
}

# --- 4. View state: Points becomes the active tab, with fresh selections ---
$templates.Range("A1:F39").Select() | Out-Null

$points.Activate() | Out-Null
$points.Range("D2").Select() | Out-Null
